$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, pushing existing rows 12-25 down to 13-26.
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with data.
$ws.Cells.Item(12, 1).Value = 5
$ws.Cells.Item(12, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(12, 3).Value = "Maule"
$ws.Cells.Item(12, 4).Value = 44741
$ws.Cells.Item(12, 5).Value = 7
$ws.Cells.Item(12, 6).Value = 100112040
$ws.Cells.Item(12, 7).Value = "Cilantro"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 150
$ws.Cells.Item(12, 11).Value = 9000
$ws.Cells.Item(12, 12).Value = 9000
$ws.Cells.Item(12, 13).Value = 9000
$ws.Cells.Item(12, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(12, 15).Value = "Región Metropolitana"
$ws.Cells.Item(12, 16).Value = 250
$ws.Cells.Item(12, 17).Value = 36
$ws.Cells.Item(12, 18).Value = "Hortaliza"
